$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt 10; $i++) {
  $ws.Cells.Item(1, $i+1).Value = $headers[$i] + "_FV2310"
  $ws.Cells.Item(1, $i+12).Value = $headers[$i] + "_FV2404"
}

$headerRange = $ws.Range("A1:U1")
$headerRange.ClearFormats()

$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U58"), 0, 1)
$tbl.Name = "Table1"

# restore header formatting (match original s="1": bold font, D9D9D9 fill, thin border all sides, center+wrap)
$headerRange.Font.Bold = $true
$headerRange.Interior.Color = 14277081
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.WrapText = $true
